$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: add "S" / "Summary" in columns A/B, matching the style already
# used by the existing C16/D16 cells (style index 4) ---
$ws.Range("C16:D16").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$ws.Range("A16").Value = "S"
$ws.Range("B16").Value = "Summary"

# --- New row 24: "Mental Health" header, merged across A24:B24, centered ---
$ws.Range("A24:B24").Merge()
$ws.Range("A24:B24").HorizontalAlignment = -4108
$ws.Range("A24").Value = "Mental Health"
$ws.Range("A24").Font.Name = "Arial"

# --- View: move selection near the bottom of the newly extended table ---
$ws.Range("A25").Select()
